# fix(module3): use uncon_planned_qty for future production; keep produced for today
# Update the ProductionPlan sheet: revise row 2 quantities and add row 3 for MAT_B.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductionPlan")

# --- Update existing row 2 (MAT_A) quantities ---
$ws.Range("G2").Value = 710
$ws.Range("H2").Value = 710
$ws.Range("J2").Value = 675

# --- Add new row 3 (MAT_B) ---
$ws.Range("A3").Value = "MAT_B"
$ws.Range("B3").Value = "PLANT_001"
$ws.Range("C3").Value = "LINE_B"

$ws.Range("D3").Value = 45293
$ws.Range("E3").Value = 45294
$ws.Range("F3").Value = 45295
$ws.Range("D3").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E3").NumberFormat = $ws.Range("E2").NumberFormat
$ws.Range("F3").NumberFormat = $ws.Range("F2").NumberFormat

$ws.Range("G3").Value = 104
$ws.Range("H3").Value = 104
# I3 stays empty (mirrors I2's blank changeover_id) but must still exist as
# a real cell in sheetData so the used-range dimension covers J3. A bare
# Formula = "'" materialises an empty text cell; resetting the style back
# to Normal afterwards drops the quote-prefix flag that Formula assignment
# implicitly sets.
$ws.Range("I3").Formula = "'"
$ws.Range("I3").Style = "Normal"
$ws.Range("J3").Value = 92
